$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking text value into a cell while keeping it
# stored as text (matching the workbook's inlineStr cells), and without
# leaving a visible style change on the cell itself.
function Set-TextValue($range, [string]$value) {
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "60.361.21"
$ws.Range("E2").Value = "  -3.02%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.305.54"
$ws.Range("E3").Value = "  -3.53%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "555.75"
$ws.Range("E5").Value = "  -4.23%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "142.38"
$ws.Range("E6").Value = "  -7.58%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").Value = "3.304.34"
$ws.Range("E8").Value = "  -3.50%  "

# Row 9 - XRP
Set-TextValue $ws.Range("D9") "0.469"
$ws.Range("E9").Value = "  -3.44%  "

# Row 10 - Toncoin
Set-TextValue $ws.Range("D10") "7.89"
$ws.Range("E10").Value = "  -2.60%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -5.16%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("D12") "0.411"
$ws.Range("E12").Value = "  -2.50%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.872.35"
$ws.Range("E13").Value = "  -3.50%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.27%  "

# Row 15 - Avalanche
Set-TextValue $ws.Range("D15") "26.56"
$ws.Range("E15").Value = "  -8.10%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.305.53"
$ws.Range("E16").Value = "  -3.54%  "

# Row 17 - ShibaInu
$ws.Range("E17").Value = "  -5.50%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "60.320.34"
$ws.Range("E18").Value = "  -3.03%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "6.11"
$ws.Range("E19").Value = "  -6.13%  "

# Row 20 - Chainlink
Set-TextValue $ws.Range("D20") "13.62"
$ws.Range("E20").Value = "  -5.95%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "8.55"
$ws.Range("E21").Value = "  -5.22%  "

# Row 22 - BitcoinCash
Set-TextValue $ws.Range("D22") "373.44"
$ws.Range("E22").Value = "  -2.38%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  +0.02%  "

# Row 24 - Litecoin
Set-TextValue $ws.Range("D24") "72.10"
$ws.Range("E24").Value = "  -5.19%  "

# Row 25 - Polygon
Set-TextValue $ws.Range("D25") "0.533"
$ws.Range("E25").Value = "  -6.93%  "

# Row 26 - WrappedeETH
$ws.Range("D26").Value = "3.442.92"
$ws.Range("E26").Value = "  -3.39%  "

# Row 27 - PEPE
Set-TextValue $ws.Range("D27") "0.0000103"
$ws.Range("E27").Value = "  -8.79%  "

# Row 28 - Kaspa
Set-TextValue $ws.Range("D28") "0.175"
$ws.Range("E28").Value = "  -1.86%  "

# Row 29 - Binance-PegBSC-USD
$ws.Range("E29").Value = "  +0.31%  "

# Row 30 - RenderToken
Set-TextValue $ws.Range("D30") "7.09"
$ws.Range("E30").Value = "  -7.89%  "

# Row 31 - USDe
Set-TextValue $ws.Range("D31") "1.00"
$ws.Range("E31").Value = "  -0.02%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  -5.50%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "7.25"
$ws.Range("E33").Value = "  -7.92%  "

# Row 34 - EthereumClassic
$ws.Range("E34").Value = "  -3.06%  "

# Row 35 - Fetch.AI
$ws.Range("E35").Value = "  -5.65%  "

# Row 36 - NEARProtocol
Set-TextValue $ws.Range("D36") "5.06"
$ws.Range("E36").Value = "  -8.76%  "

# Row 37 - Monero
Set-TextValue $ws.Range("D37") "165.94"
$ws.Range("E37").Value = "  -1.47%  "

# Row 38 - ImmutableX
Set-TextValue $ws.Range("D38") "1.51"
$ws.Range("E38").Value = "  -5.66%  "

# Row 39 - Aptos
Set-TextValue $ws.Range("D39") "6.62"
$ws.Range("E39").Value = "  -5.23%  "

# Row 40 - was RenzoRestakedETH, now Hedera
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D40") "0.0722"
$ws.Range("E40").Value = "  -7.73%  "

# Row 41 - was Hedera, now RenzoRestakedETH
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.335.02"
$ws.Range("E41").Value = "  -3.72%  "

# Row 42 - EnergySwap
Set-TextValue $ws.Range("D42") "25.26"
$ws.Range("E42").Value = "  -18.58%  "

# Row 43 - OKB
Set-TextValue $ws.Range("D43") "41.85"
$ws.Range("E43").Value = "  -2.25%  "

# Row 44 - Mantle
Set-TextValue $ws.Range("D44") "0.748"
$ws.Range("E44").Value = "  -4.36%  "

# Row 45 - ONDO
Set-TextValue $ws.Range("D45") "1.13"
$ws.Range("E45").Value = "  -3.45%  "

# Row 46 - Filecoin
Set-TextValue $ws.Range("D46") "4.12"
$ws.Range("E46").Value = "  -6.91%  "

# Row 47 - Stacks
$ws.Range("E47").Value = "  -6.71%  "

# Row 48 - FirstDigitalUSD
$ws.Range("E48").Value = "  -0.04%  "

# Row 49 - Maker
$ws.Range("D49").Value = "2.322.46"
$ws.Range("E49").Value = "  -9.17%  "

# Row 50 - Cosmos
$ws.Range("E50").Value = "  -6.87%  "

# Row 51 - InjectiveProtocol
Set-TextValue $ws.Range("D51") "21.58"
$ws.Range("E51").Value = "  -8.05%  "
